# Updated symbol list on Sun Feb  5 22:45:29 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row on Sheet1. These columns hold plain text (e.g. "328.05",
# "-1.08%") rather than real numbers/percentages, so each new value is
# written with a leading apostrophe to force Excel to keep it as literal
# text instead of auto-converting it to a Number or Percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.18"
$ws.Range("E2").Value = "'-0.92%"

$ws.Range("D3").Value = "'43.81"
$ws.Range("E3").Value = "'5.99%"

$ws.Range("D4").Value = "'5.571"
$ws.Range("E4").Value = "'-2.13%"

$ws.Range("D5").Value = "'0.08108"
$ws.Range("E5").Value = "'-3.75%"

$ws.Range("D6").Value = "'8.663"
$ws.Range("E6").Value = "'-1.74%"

$ws.Range("D7").Value = "'4.281"
$ws.Range("E7").Value = "'-4.84%"

$ws.Range("D8").Value = "'1.895"
$ws.Range("E8").Value = "'-4.59%"

$ws.Range("E9").Value = "'-3.73%"

$ws.Range("D10").Value = "'0.9383"
$ws.Range("E10").Value = "'1.05%"

$ws.Range("D11").Value = "'0.1173"
$ws.Range("E11").Value = "'-6.26%"

$ws.Range("D12").Value = "'0.1891"
$ws.Range("E12").Value = "'-4.62%"

$ws.Range("D13").Value = "'0.09606"
$ws.Range("E13").Value = "'2.10%"

$ws.Range("D14").Value = "'0.04200"
$ws.Range("E14").Value = "'6.09%"

$ws.Range("D15").Value = "'0.1067"
$ws.Range("E15").Value = "'0.43%"

$ws.Range("D16").Value = "'0.001268"
$ws.Range("E16").Value = "'-3.63%"

$ws.Range("D17").Value = "'0.005927"
$ws.Range("E17").Value = "'-2.93%"

$ws.Range("D18").Value = "'3.568"
$ws.Range("E18").Value = "'3.87%"

$ws.Range("D19").Value = "'0.3485"
$ws.Range("E19").Value = "'-0.76%"

$ws.Range("D20").Value = "'8.564"
$ws.Range("E20").Value = "'-6.52%"

$ws.Range("D21").Value = "'0.1361"
$ws.Range("E21").Value = "'-0.11%"

$ws.Range("E22").Value = "'3.15%"

$ws.Range("D23").Value = "'0.04324"
$ws.Range("E23").Value = "'-2.25%"

$ws.Range("E24").Value = "'-0.52%"

$ws.Range("D25").Value = "'0.004396"
$ws.Range("E25").Value = "'0.47%"

$ws.Range("D26").Value = "'0.0001230"
$ws.Range("E26").Value = "'3.15%"

$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-0.02%"

$ws.Range("D39").Value = "'0.02669"
$ws.Range("E39").Value = "'-5.79%"

$ws.Range("D40").Value = "'0.05476"
$ws.Range("E40").Value = "'-1.02%"

$ws.Range("D41").Value = "'0.01144"
$ws.Range("E41").Value = "'27.80%"

$ws.Range("D42").Value = "'0.007684"
$ws.Range("E42").Value = "'-2.73%"

$ws.Range("D43").Value = "'0.1398"
$ws.Range("E43").Value = "'-2.99%"

$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'1.12%"

$ws.Range("D45").Value = "'0.009639"
$ws.Range("E45").Value = "'-6.36%"

$ws.Range("D46").Value = "'0.00007057"
$ws.Range("E46").Value = "'-3.35%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.02%"

$ws.Range("D48").Value = "'0.003559"
$ws.Range("E48").Value = "'10.90%"

$ws.Range("D49").Value = "'0.002270"
$ws.Range("E49").Value = "'-0.35%"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.02%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.02%"
